# Updating price policy code
#
# The "Electric_boiler" technology line is dropped from each of the four
# per-technology cost sheets (its row is deleted, shifting the remaining
# technologies - Gas_CHP, Gas_boiler, Grid, Heat_pump, Solar_PV,
# Solar_thermal - up by one row), and the recomputed cost figures are
# written in for the new policy run.

$wb = $excel.ActiveWorkbook

# --- Operating_cost_per_technology -----------------------------------
$ws = $wb.Worksheets.Item("Operating_cost_per_technology")
$ws.Rows("1:1").Delete()
$ws.Range("B1").Value = 67120.151889195142
$ws.Range("B2").Value = 88143.333816788523

# --- Maintenance_cost_per_technology -----------------------------------
$ws = $wb.Worksheets.Item("Maintenance_cost_per_technology")
$ws.Rows("1:1").Delete()
$ws.Range("B1").Value = 12826.661026025011
$ws.Range("B2").Value = 9206.0815319759586
$ws.Range("B5").Value = 29319.134086180209

# --- Capital_cost_per_technology ---------------------------------------
$ws = $wb.Worksheets.Item("Capital_cost_per_technology")
$ws.Rows("1:1").Delete()
$ws.Range("B1").Value = 20197.873257037721
$ws.Range("B2").Value = 29773.839895816691
$ws.Range("B5").Value = 71214.076980455284

# --- Total_cost_per_technology ------------------------------------------
$ws = $wb.Worksheets.Item("Total_cost_per_technology")
$ws.Rows("1:1").Delete()
$ws.Range("B1").Value = 100144.68617225788
$ws.Range("B2").Value = 127123.25524458117
$ws.Range("B5").Value = 100533.2110666355

# --- Operating_cost_grid --------------------------------------------------
$ws = $wb.Worksheets.Item("Operating_cost_grid")
$ws.Range("A1").Value = 169329.0430287901

# --- Total_cost_grid -------------------------------------------------------
$ws = $wb.Worksheets.Item("Total_cost_grid")
$ws.Range("A1").Value = 169329.0430287901

# --- Capital_cost_per_storage -----------------------------------------------
$ws = $wb.Worksheets.Item("Capital_cost_per_storage")
$ws.Range("B2").Value = 3908.5031409431786

# --- Total_cost_per_storage --------------------------------------------------
$ws = $wb.Worksheets.Item("Total_cost_per_storage")
$ws.Range("B2").Value = 3908.5031409431786

# --- Income_via_exports ---------------------------------------------------
$ws = $wb.Worksheets.Item("Income_via_exports")
$ws.Range("A1").Value = 45651.636692740649
